$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Rank" column (column B). This shifts Tuesday..Monday
# (previously C:I) left into B:H, preserving each cell's contents and
# the header style.
$ws.Range("B1").EntireColumn.Delete()
